$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / non-numeric-looking updates ---
$ws.Range("D2").Value = "66.891.63"
$ws.Range("E2").Value = "  +2.99%  "
$ws.Range("D3").Value = "3.202.23"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +4.31%  "
$ws.Range("E6").Value = "  +5.90%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +6.42%  "
$ws.Range("D9").Value = "3.204.51"
$ws.Range("E9").Value = "  +1.93%  "
$ws.Range("E10").Value = "  +2.14%  "
$ws.Range("E11").Value = "  -3.59%  "
$ws.Range("E12").Value = "  +3.82%  "
$ws.Range("E13").Value = "  +2.49%  "
$ws.Range("E14").Value = "  +6.17%  "
$ws.Range("D15").Value = "3.739.70"
$ws.Range("E15").Value = "  +2.06%  "
$ws.Range("D16").Value = "66.884.40"
$ws.Range("E16").Value = "  +3.05%  "
$ws.Range("E17").Value = "  +5.28%  "
$ws.Range("D18").Value = "3.214.40"
$ws.Range("E18").Value = "  +1.98%  "
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("E20").Value = "  +4.31%  "
$ws.Range("E21").Value = "  +3.07%  "
$ws.Range("E22").Value = "  +4.58%  "
$ws.Range("E23").Value = "  +6.88%  "
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("E25").Value = "  +1.66%  "
$ws.Range("E27").Value = "  +2.73%  "
$ws.Range("E28").Value = "  +4.03%  "
$ws.Range("E29").Value = "  +10.49%  "
$ws.Range("E30").Value = "  +8.28%  "
$ws.Range("E31").Value = "  +9.93%  "
$ws.Range("E32").Value = "  +3.32%  "
$ws.Range("E33").Value = "  +3.32%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("E35").Value = "  +1.87%  "
$ws.Range("E36").Value = "  +10.93%  "
$ws.Range("E37").Value = "  +0.52%  "
$ws.Range("E38").Value = "  +2.66%  "
$ws.Range("E39").Value = "  +3.99%  "
$ws.Range("E40").Value = "  +9.45%  "
$ws.Range("E41").Value = "  +2.37%  "
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("D43").Value = "0.0₃0688"
$ws.Range("E43").Value = "  +16.54%  "
$ws.Range("E44").Value = "  +7.76%  "
$ws.Range("E45").Value = "  +3.01%  "
$ws.Range("D46").Value = "2.905.88"
$ws.Range("E46").Value = "  -2.97%  "
$ws.Range("E47").Value = "  +1.92%  "
$ws.Range("E48").Value = "  +10.91%  "
$ws.Range("E49").Value = "  +4.01%  "
$ws.Range("B50").Value = "USDe"
$ws.Range("C50").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("E51").Value = "  +4.84%  "

# --- Numeric-looking price values: force Text so Excel keeps the literal
#     string instead of auto-converting to a Number, then restore the
#     cell style so no stray number-format is left behind. ---
$numericCells = @("D5", "D6", "D10", "D11", "D12", "D14", "D17", "D20", "D21", "D22", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D35", "D36", "D37", "D38", "D41", "D42", "D47", "D48", "D49", "D50", "D51")
foreach ($c in $numericCells) { $ws.Range($c).NumberFormat = "@" }
$ws.Range("D5").Value = "604.34"
$ws.Range("D6").Value = "157.29"
$ws.Range("D10").Value = "0.161"
$ws.Range("D11").Value = "5.91"
$ws.Range("D12").Value = "0.517"
$ws.Range("D14").Value = "39.31"
$ws.Range("D17").Value = "7.50"
$ws.Range("D20").Value = "523.91"
$ws.Range("D21").Value = "15.49"
$ws.Range("D22").Value = "0.745"
$ws.Range("D24").Value = "15.10"
$ws.Range("D25").Value = "85.47"
$ws.Range("D27").Value = "9.29"
$ws.Range("D28").Value = "3.03"
$ws.Range("D29").Value = "2.40"
$ws.Range("D30").Value = "3.01"
$ws.Range("D31").Value = "7.01"
$ws.Range("D32").Value = "28.36"
$ws.Range("D35").Value = "6.59"
$ws.Range("D36").Value = "526.08"
$ws.Range("D37").Value = "55.16"
$ws.Range("D38").Value = "0.0907"
$ws.Range("D41").Value = "8.93"
$ws.Range("D42").Value = "2.91"
$ws.Range("D47").Value = "28.74"
$ws.Range("D48").Value = "2.76"
$ws.Range("D49").Value = "0.118"
$ws.Range("D50").Value = "0.999"
$ws.Range("D51").Value = "2.35"
foreach ($c in $numericCells) { $ws.Range($c).Style = "Normal" }
